$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-gender")

# Update is_active (column D) from TRUE to FALSE for the "Others" gender rows
# Row 4: OTH / eng, Row 7: OTH / ara, Row 10: OTH / fra
$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the active cell selection to D12 (matches saved sheetView selection)
$ws.Range("D12").Select()

$wb.Save()
